$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new cell C23 = "state" (no special style, matching the diff exactly)
$ws.Range("C23").Value = "state"

# Update sheet view: scroll so A10 is top-left, and change selection to E19
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E19").Select()
